$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.331892967224121
$ws.Range("B1").Value = 2.447057962417603
$ws.Range("C1").Value = 2.03644061088562
$ws.Range("D1").Value = 2.096269607543945
$ws.Range("E1").Value = 2.416609048843384
